# Fruta / hortaliza, semanal
# Insert a new weekly data row for "Vega Monumental Concepción - Kiwi" at
# row 321 (pushing the existing rows 321:350 down to 322:351), then
# populate the new row with this week's reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 321 - this shifts rows 321:350 down
# to 322:351 and grows the sheet's used range to A1:T351.
$ws.Rows("321").Insert()

# Populate the newly inserted row 321 with the new record's values.
$ws.Range("A321").Value = 11
$ws.Range("B321").Value = "Vega Monumental Concepción"
$ws.Range("C321").Value = "Bíobío"
$ws.Range("D321").Value = 45134
$ws.Range("E321").Value = 8
$ws.Range("F321").Value = "Fruta"
$ws.Range("G321").Value = 100101
$ws.Range("H321").Value = "Berries"
$ws.Range("I321").Value = 100101007
$ws.Range("J321").Value = "Kiwi"
$ws.Range("K321").Value = "Hayward"
$ws.Range("L321").Value = "Primera"
$ws.Range("M321").Value = 270
$ws.Range("N321").Value = 11000
$ws.Range("O321").Value = 12000
$ws.Range("P321").Value = 11444
$ws.Range("Q321").Value = "$/bandeja 10 kilos"
$ws.Range("R321").Value = "Provincia de Curicó"
$ws.Range("S321").Value = 1144
$ws.Range("T321").Value = 10
